$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Servicio de consultoría estadística."
$ws.Range("B2").Value = "servicio"
$ws.Range("D2").Value = 33000

$ws.Range("D3").Select()
